# "Generate Report for Handback" - refresh the handoff/handback timestamps
# recorded for the zh-cn and de-de language reports.

$wb = $excel.ActiveWorkbook

$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("E2").Value = "2016-03-21 02:54:43"
$zh.Range("H2").Value = "2016-03-21 02:55:03"

$de = $wb.Worksheets.Item("de-de")
$de.Range("E2").Value = "2016-03-21 02:54:46"
$de.Range("H2").Value = "2016-03-21 02:55:09"
